$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Insert two new rows (9 & 10) after "Skip level button" (row 8) for
#    the new "Opening screen" / "Closing screen" tasks. Formatting is
#    copied from an existing "light" row (row 14 - "Make strings look
#    like visual strings") so the new rows pick up matching borders /
#    wrap-text formatting.
# ---------------------------------------------------------------------
$ws.Range("B14:E14").Copy()
$ws.Rows("9").Insert()
$ws.Range("B9:E9").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B15:E15").Copy()
$ws.Rows("10").Insert()
$ws.Range("B10:E10").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B9").Value = "Opening screen"
$ws.Range("C9").Value = "Dylan"
$ws.Range("D9").Value = "High"
$ws.Range("E9").Value = "Shows on startup."

$ws.Range("B10").Value = "Closing screen"
$ws.Range("C10").Value = "Dylan"
$ws.Range("D10").Value = "High"
$ws.Range("E10").Value = "Show on gameover or when you complete"

# ---------------------------------------------------------------------
# 2) Insert a new row 7 (shifting everything below it down) for the
#    "Fix boxes moving through each other when grabbed" task, copying
#    formatting from the neighbouring "normal" row (row 8, "Fix IK for
#    climbing") and then giving it a taller row height for the longer
#    comment text.
# ---------------------------------------------------------------------
$ws.Range("B8:E8").Copy()
$ws.Rows("7").Insert()
$ws.Range("B7:E7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B7").Value = "Fix boxes moving through each other when grabbed"
$ws.Range("C7").Value = "Tim"
$ws.Range("D7").Value = "Mid"
$ws.Range("E7").Value = "Hard cause I know to move them we disabled the collider, but has to be fixed"
$ws.Rows("7").RowHeight = 30

# ---------------------------------------------------------------------
# 3) Mark completed tasks: highlight "Skip level button" (now row 9)
#    and "Add colliders to stage poles" (now row 15) in green to show
#    they've been ticked off.
# ---------------------------------------------------------------------
$ws.Range("B9:E9").Interior.Color = 5296274
$ws.Range("B15:E15").Interior.Color = 5296274

# ---------------------------------------------------------------------
# 4) Refresh the active selection to match where the user left off.
# ---------------------------------------------------------------------
$ws.Range("I9").Select()

Write-Host "Cosmetic changes applied"
